$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = "No"
$ws.Range("C22").Value = "Yes"
$ws.Range("C56").Value = "Yes"
$ws.Range("C57").Value = "No"
$ws.Range("C65").Value = "Yes"
$ws.Range("C67").Value = "No"
$ws.Range("C76").Value = "No"
$ws.Range("C83").Value = "Yes"

$ws.Range("C84").Select()
